$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the source citation text (shared string used by the last row) ---
# Original: Federal State Statistic Service (FSSS), "Results of SME Census in Russia", 2010. ...
# New:      European Investment Bank - EIB, Small and Medium Entrepreneurship in Russia (2013) ...
$ws.Range("A46").Value2 = 'European Investment Bank - EIB, Small and Medium Entrepreneurship in Russia (2013) available at http://www.eib.org/attachments/efs/econ_study_small_and_medium_entrepreneurship_in_russia_en.pdf'

# --- Insert a new (blank) row before row 40, shifting the source block down by one ---
$ws.Rows("40").Insert()

# Remember the hyperlinked URL text currently sitting in A42 (it carries the
# "HyperLink" style which must not survive in the new layout)
$url = $ws.Range("A42").Value2

# Remove that hyperlink-styled row entirely; rows below shift back up
$ws.Rows("42").Delete()

# Insert a fresh row at 43 (it inherits the plain "source" style from the row
# above it, i.e. the same style used by the rest of the block)
$ws.Rows("43").Insert()

# Put the URL text back, now on a plain (non-hyperlink) styled cell
$ws.Range("A43").Value2 = $url

# Drop the hyperlink itself - the new text will just be plain text
$ws.Hyperlinks.Delete()
